$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh inserts two new daily price records (rows 38 and 41 in
# the final layout) ahead of the previously-reported data, pushing the
# existing rows down accordingly.
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(41).Insert()

# New row 38: Arándano (blue), Primera, Región de Ñuble
$ws.Cells.Item(38, 1).Value = 11
$ws.Cells.Item(38, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(38, 3).Value = "Bíobío"
$ws.Cells.Item(38, 4).Value = 44518
$ws.Cells.Item(38, 5).Value = 8
$ws.Cells.Item(38, 6).Value = "Fruta"
$ws.Cells.Item(38, 7).Value = 100101
$ws.Cells.Item(38, 8).Value = "Berries"
$ws.Cells.Item(38, 9).Value = 100101001
$ws.Cells.Item(38, 10).Value = "Arándano (blue)"
$ws.Cells.Item(38, 11).Value = "Sin especificar"
$ws.Cells.Item(38, 12).Value = "Primera"
$ws.Cells.Item(38, 13).Value = 200
$ws.Cells.Item(38, 14).Value = 7000
$ws.Cells.Item(38, 15).Value = 7500
$ws.Cells.Item(38, 16).Value = 7250
$ws.Cells.Item(38, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(38, 18).Value = "Región de Ñuble"
$ws.Cells.Item(38, 19).Value = 3625
$ws.Cells.Item(38, 20).Value = 2

# New row 41: Arándano (blue), Primera, Provincia de Linares
$ws.Cells.Item(41, 1).Value = 11
$ws.Cells.Item(41, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(41, 3).Value = "Bíobío"
$ws.Cells.Item(41, 4).Value = 44523
$ws.Cells.Item(41, 5).Value = 8
$ws.Cells.Item(41, 6).Value = "Fruta"
$ws.Cells.Item(41, 7).Value = 100101
$ws.Cells.Item(41, 8).Value = "Berries"
$ws.Cells.Item(41, 9).Value = 100101001
$ws.Cells.Item(41, 10).Value = "Arándano (blue)"
$ws.Cells.Item(41, 11).Value = "Sin especificar"
$ws.Cells.Item(41, 12).Value = "Primera"
$ws.Cells.Item(41, 13).Value = 150
$ws.Cells.Item(41, 14).Value = 5000
$ws.Cells.Item(41, 15).Value = 5500
$ws.Cells.Item(41, 16).Value = 5267
$ws.Cells.Item(41, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(41, 18).Value = "Provincia de Linares"
$ws.Cells.Item(41, 19).Value = 2634
$ws.Cells.Item(41, 20).Value = 2

# Match the date-format style used by the rest of column D ("Fecha") on the
# two newly inserted rows.
$ws.Range("D38").NumberFormat = $ws.Range("D39").NumberFormat
$ws.Range("D41").NumberFormat = $ws.Range("D42").NumberFormat
